$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the 3rd column of "Tabla3" (H6:J8) from "X * F" to "X + F".
#    Writing the header cell directly also renames the table column (the
#    ListColumns.Item(n).Name setter is a no-op in this host, so we go
#    through the worksheet cell instead).
# ---------------------------------------------------------------------------
$ws.Range("J6").Value = "X + F"

# ---------------------------------------------------------------------------
# 2. Move the existing "Tabla6" (idempotence law, X / X+X / X*X) from
#    H17:J19 down to H26:J28, and relocate its "Ley de idempotencia" caption
#    from H16 to H25.
# ---------------------------------------------------------------------------
$tabla6 = $ws.ListObjects.Item("Tabla6")
$tabla6.Resize($ws.Range("H26:J28"))

$ws.Range("H26").Value = "X"
$ws.Range("I26").Value = "X + X"
$ws.Range("J26").Value = "X * X"
$ws.Range("H27").Value = "t"
$ws.Range("I27").Value = "t"
$ws.Range("J27").Value = "t"
$ws.Range("H28").Value = "f"
$ws.Range("I28").Value = "f"
$ws.Range("J28").Value = "f"

$ws.Range("H25").Value = "Ley de idempotencia"

# Clear the now-vacated old locations.
$ws.Range("H16:J19").ClearContents()

# ---------------------------------------------------------------------------
# 3. The old "Tabla5" (H12:J14 -- X / X+'X / X*'X) is replaced by two new,
#    smaller complement-law tables ("Tabla5" reused name at H16:J18 and
#    "Tabla59" at H20:J22). Remove the old table/data first.
# ---------------------------------------------------------------------------
$oldTabla5 = $ws.ListObjects.Item("Tabla5")
$oldTabla5.Delete()
$ws.Range("H11:J14").ClearContents()

$ws.Range("H15").Value = "Ley de complementos"

# ---------------------------------------------------------------------------
# 4. New table: Tabla311 (identity law, AND) at H10:J12.
# ---------------------------------------------------------------------------
$ws.Range("H10").Value = "X"
$ws.Range("I10").Value = "F"
$ws.Range("J10").Value = "X * F"
$ws.Range("H11").Value = "t"
$ws.Range("I11").Value = "f"
$ws.Range("J11").Value = "f"
$ws.Range("H12").Value = "f"
$ws.Range("I12").Value = "f"
$ws.Range("J12").Value = "f"

$tabla311 = $ws.ListObjects.Add(1, $ws.Range("H10:J12"), $null, 1)
$tabla311.ShowAutoFilterDropDown = $false
$tabla311.TableStyle = "TableStyleDark2"
$tabla311.Name = "Tabla311"

# ---------------------------------------------------------------------------
# 5. New table: Tabla5 (reused name, complement law OR) at H16:J18.
# ---------------------------------------------------------------------------
$ws.Range("H16").Value = "X "
$ws.Range("I16").Value = "´X"
$ws.Range("J16").Value = "X + ´X"
$ws.Range("H17").Value = "t"
$ws.Range("I17").Value = "f"
$ws.Range("J17").Value = "t"
$ws.Range("H18").Value = "f"
$ws.Range("I18").Value = "t"
$ws.Range("J18").Value = "t"

$tabla5new = $ws.ListObjects.Add(1, $ws.Range("H16:J18"), $null, 1)
$tabla5new.ShowAutoFilterDropDown = $false
$tabla5new.TableStyle = "TableStyleDark2"
$tabla5new.Name = "Tabla5"

# ---------------------------------------------------------------------------
# 6. New table: Tabla59 (complement law AND) at H20:J22.
# ---------------------------------------------------------------------------
$ws.Range("H20").Value = "X "
$ws.Range("I20").Value = "´X"
$ws.Range("J20").Value = "X * ´X"
$ws.Range("H21").Value = "t"
$ws.Range("I21").Value = "f"
$ws.Range("J21").Value = "f"
$ws.Range("H22").Value = "f"
$ws.Range("I22").Value = "t"
$ws.Range("J22").Value = "f"

$tabla59 = $ws.ListObjects.Add(1, $ws.Range("H20:J22"), $null, 1)
$tabla59.ShowAutoFilterDropDown = $false
$tabla59.TableStyle = "TableStyleDark2"
$tabla59.Name = "Tabla59"

# ---------------------------------------------------------------------------
# 7. New table: Tabla13 (associative law, AND) at B16:F24 -- mirrors Tabla1
#    (B5:F13, the OR version) one-for-one.
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "X"
$ws.Range("C16").Value = "Y"
$ws.Range("D16").Value = "Z"
$ws.Range("E16").Value = "(X * Y) * Z"
$ws.Range("F16").Value = "X * (Y * Z)"

$ws.Range("B17").Value = "t"
$ws.Range("C17").Value = "t"
$ws.Range("D17").Value = "t"
$ws.Range("E17").Value = "t"
$ws.Range("F17").Value = "t"

$ws.Range("B18").Value = "t"
$ws.Range("C18").Value = "t"
$ws.Range("D18").Value = "f"
$ws.Range("E18").Value = "f"
$ws.Range("F18").Value = "f"

$ws.Range("B19").Value = "t"
$ws.Range("C19").Value = "f"
$ws.Range("D19").Value = "t"
$ws.Range("E19").Value = "f"
$ws.Range("F19").Value = "f"

$ws.Range("B20").Value = "t"
$ws.Range("C20").Value = "f"
$ws.Range("D20").Value = "f"
$ws.Range("E20").Value = "f"
$ws.Range("F20").Value = "f"

$ws.Range("B21").Value = "f"
$ws.Range("C21").Value = "t"
$ws.Range("D21").Value = "t"
$ws.Range("E21").Value = "f"
$ws.Range("F21").Value = "f"

$ws.Range("B22").Value = "f"
$ws.Range("C22").Value = "t"
$ws.Range("D22").Value = "f"
$ws.Range("E22").Value = "f"
$ws.Range("F22").Value = "f"

$ws.Range("B23").Value = "f"
$ws.Range("C23").Value = "f"
$ws.Range("D23").Value = "t"
$ws.Range("E23").Value = "f"
$ws.Range("F23").Value = "f"

$ws.Range("B24").Value = "f"
$ws.Range("C24").Value = "f"
$ws.Range("D24").Value = "f"
$ws.Range("E24").Value = "f"
$ws.Range("F24").Value = "f"

$tabla13 = $ws.ListObjects.Add(1, $ws.Range("B16:F24"), $null, 1)
$tabla13.ShowAutoFilterDropDown = $false
$tabla13.TableStyle = "TableStyleDark2"
$tabla13.Name = "Tabla13"

# ---------------------------------------------------------------------------
# 8. New table: Tabla712 (De Morgan's law, AND) at M14:P18 -- mirrors Tabla7
#    (M7:P11, the OR version) one-for-one.
# ---------------------------------------------------------------------------
$ws.Range("M14").Value = "X"
$ws.Range("N14").Value = "Y"
$ws.Range("O14").Value = "( X * Y)^2"
$ws.Range("P14").Value = "´X * 'Y"

$ws.Range("M15").Value = "t"
$ws.Range("N15").Value = "t"
$ws.Range("O15").Value = "t"
$ws.Range("P15").Value = "t"

$ws.Range("M16").Value = "t"
$ws.Range("N16").Value = "f"
$ws.Range("O16").Value = "f"
$ws.Range("P16").Value = "f"

$ws.Range("M17").Value = "f"
$ws.Range("N17").Value = "t"
$ws.Range("O17").Value = "f"
$ws.Range("P17").Value = "f"

$ws.Range("M18").Value = "f"
$ws.Range("N18").Value = "f"
$ws.Range("O18").Value = "f"
$ws.Range("P18").Value = "f"

$tabla712 = $ws.ListObjects.Add(1, $ws.Range("M14:P18"), $null, 1)
$tabla712.ShowAutoFilterDropDown = $false
$tabla712.TableStyle = "TableStyleDark2"
$tabla712.Name = "Tabla712"

# ---------------------------------------------------------------------------
# 9. View state: scroll so row 4 is at the top, and select R17 (matches the
#    saved workbook view in the edited file).
# ---------------------------------------------------------------------------
$ws.Range("R17").Select()
$excel.ActiveWindow.ScrollRow = 4
